{"js": "// Updated Red Bull end date: \" - Present\" -> \" \\u2013 Sep 2025\"\nconst body = context.document.body;\nconst results = body.search(\" - Present\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find ' - Present' text to replace.\");\n}\n\n// Replace the matched text (\" - Present\") with \" \\u2013 Sep 2025\" (en dash).\nresults.items[0].insertText(\" \\u2013 Sep 2025\", \"Replace\");\nawait context.sync();\n", "ps1": "# Updated Red Bull end Date\n# Change the Oracle Red Bull Racing end date from \" - Present\" to \" \\u2013 Sep 2025\" (en dash).\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Text = \" - Present\"\n$range.Find.MatchCase = $true\n$range.Find.MatchWholeWord = $false\n$range.Find.Forward = $true\n\n$found = $range.Find.Execute()\nif ($found) {\n    # $range now spans exactly the matched text (\" - Present\"); overwrite it\n    # with \" \\u2013 Sep 2025\" (space, en dash, space, new end date).\n    $range.Text = \" \" + [char]0x2013 + \" Sep 2025\"\n}\n"}
